# Update "想去人数" (interested-count) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new value for column F, per-sheet (F13 differs between the two sheets).
$commonUpdates = @{
    2  = 1550
    3  = 44
    4  = 1011
    5  = 12
    7  = 2553
    9  = 1598
    11 = 182
    15 = 49
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $commonUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $commonUpdates[$row]
    }

    if ($sheetName -eq "展览") {
        $ws.Cells.Item(13, 6).Value = 505
    } elseif ($sheetName -eq "全部类型") {
        $ws.Cells.Item(13, 6).Value = 506
    }
}
